$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "73.022.39"
$ws.Cells.Item(2, 5).Value = "  +2.92%  "

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "3.985.52"
$ws.Cells.Item(3, 5).Value = "  +1.08%  "

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.00"
$ws.Cells.Item(4, 5).Value = "  -0.05%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "609.86"
$ws.Cells.Item(5, 5).Value = "  +13.72%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "164.17"
$ws.Cells.Item(6, 5).Value = "  +11.13%  "

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.685"
$ws.Cells.Item(7, 5).Value = "  -0.16%  "

$ws.Cells.Item(8, 5).Value = "  -0.07%  "

$ws.Cells.Item(9, 5).Value = "  +2.26%  "

$ws.Cells.Item(10, 5).Value = "  +2.14%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "54.58"
$ws.Cells.Item(11, 5).Value = "  -1.12%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.0000320"
$ws.Cells.Item(12, 5).Value = "  +1.02%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "11.03"
$ws.Cells.Item(13, 5).Value = "  +4.10%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "4.629.62"
$ws.Cells.Item(14, 5).Value = "  +1.14%  "

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "3.995.31"
$ws.Cells.Item(15, 5).Value = "  +1.11%  "

$ws.Cells.Item(16, 5).Value = "  +9.34%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "14.13"
$ws.Cells.Item(17, 5).Value = "  +1.93%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "20.56"
$ws.Cells.Item(18, 5).Value = "  +0.16%  "

$ws.Cells.Item(19, 5).Value = "  +0.48%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "72.763.50"
$ws.Cells.Item(20, 5).Value = "  +2.60%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "440.31"
$ws.Cells.Item(21, 5).Value = "  +4.87%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "4.86"
$ws.Cells.Item(22, 5).Value = "  +14.57%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "96.87"
$ws.Cells.Item(23, 5).Value = "  -0.64%  "

$ws.Cells.Item(24, 5).Value = "  -3.72%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "14.40"
$ws.Cells.Item(25, 5).Value = "  -0.97%  "

$ws.Cells.Item(26, 5).Value = "  +13.01%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "11.41"
$ws.Cells.Item(27, 5).Value = "  +1.35%  "

$ws.Cells.Item(28, 5).Value = "  +1.08%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "10.43"
$ws.Cells.Item(29, 5).Value = "  -2.79%  "

$ws.Cells.Item(30, 5).Value = "  +0.22%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "7.85"
$ws.Cells.Item(31, 5).Value = "  +1.57%  "

$ws.Cells.Item(32, 5).Value = "  +4.61%  "

$ws.Cells.Item(33, 5).Value = "  +0.15%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "48.99"
$ws.Cells.Item(34, 5).Value = "  -3.64%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "71.14"
$ws.Cells.Item(35, 5).Value = "  +8.34%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "658.99"
$ws.Cells.Item(36, 5).Value = "  -2.89%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.0₃0897"
$ws.Cells.Item(37, 5).Value = "  +10.24%  "

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.438"
$ws.Cells.Item(38, 5).Value = "  -0.40%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "3.37"
$ws.Cells.Item(39, 5).Value = "  +0.07%  "

$ws.Cells.Item(40, 5).Value = "  -0.97%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "1.00"
$ws.Cells.Item(41, 5).Value = "  -0.03%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "3.33"
$ws.Cells.Item(42, 5).Value = "  +5.04%  "

$ws.Cells.Item(43, 5).Value = "  +0.14%  "

$ws.Cells.Item(44, 5).Value = "  +1.96%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "10.69"
$ws.Cells.Item(45, 5).Value = "  +5.26%  "

$ws.Cells.Item(46, 5).Value = "  +1.11%  "

$ws.Cells.Item(47, 5).Value = "  -0.04%  "

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "3.40"
$ws.Cells.Item(48, 5).Value = "  +1.75%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "2.912.47"
$ws.Cells.Item(49, 5).Value = "  +12.37%  "

$ws.Cells.Item(50, 5).Value = "  +2.23%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "3.41"
$ws.Cells.Item(51, 5).Value = "  +4.95%  "
